$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.741029
$ws.Range("H2").Value = 23.223087
$ws.Range("I2").Value = 0.4930486933812723
$ws.Range("J2").Value = 0.4930486933812723
$ws.Range("M2").Value = 3.759736666666667
$ws.Range("N2").Value = 11.27921
$ws.Range("O2").Value = 0.0683751702595819
$ws.Range("P2").Value = 0.06837517025958188
$ws.Range("Q2").Value = 29.10423056903
$ws.Range("R2").Value = 261.9380751212699
$ws.Range("S2").Value = 0.03371228835620888
$ws.Range("T2").Value = 0.03371228835620887

$ws.Range("G3").Value = 7.741029
$ws.Range("H3").Value = 23.223087
$ws.Range("I3").Value = 0.4930486933812723
$ws.Range("J3").Value = 0.4930486933812723
$ws.Range("O3").Value = 0.6514180024294648
$ws.Range("P3").Value = 0.6514180024294647
$ws.Range("Q3").Value = 277.279305741359
$ws.Range("R3").Value = 2495.513751672231
$ws.Range("S3").Value = 0.321180794942886
$ws.Range("T3").Value = 0.321180794942886

$ws.Range("G4").Value = 7.741029
$ws.Range("H4").Value = 23.223087
$ws.Range("I4").Value = 0.4930486933812723
$ws.Range("J4").Value = 0.4930486933812723
$ws.Range("O4").Value = 0.2802068273109533
$ws.Range("P4").Value = 0.2802068273109533
$ws.Range("Q4").Value = 119.27142671987
$ws.Range("R4").Value = 1073.44284047883
$ws.Range("S4").Value = 0.1381556100821773
$ws.Range("T4").Value = 0.1381556100821773

$ws.Range("I5").Value = 0.0194007766416684
$ws.Range("J5").Value = 0.0194007766416684
$ws.Range("M5").Value = 3.759736666666667
$ws.Range("N5").Value = 11.27921
$ws.Range("O5").Value = 0.0683751702595819
$ws.Range("P5").Value = 0.06837517025958188
$ws.Range("Q5").Value = 1.145210775684445
$ws.Range("R5").Value = 10.30689698116
$ws.Range("S5").Value = 0.001326531406042196
$ws.Range("T5").Value = 0.001326531406042196

$ws.Range("I6").Value = 0.0194007766416684
$ws.Range("J6").Value = 0.0194007766416684
$ws.Range("O6").Value = 0.6514180024294648
$ws.Range("P6").Value = 0.6514180024294647
$ws.Range("S6").Value = 0.01263801516549585
$ws.Range("T6").Value = 0.01263801516549585

$ws.Range("I7").Value = 0.0194007766416684
$ws.Range("J7").Value = 0.0194007766416684
$ws.Range("O7").Value = 0.2802068273109533
$ws.Range("P7").Value = 0.2802068273109533
$ws.Range("S7").Value = 0.005436230070130354
$ws.Range("T7").Value = 0.005436230070130353

$ws.Range("G8").Value = 7.654706000000001
$ws.Range("I8").Value = 0.4875505299770593
$ws.Range("J8").Value = 0.4875505299770593
$ws.Range("M8").Value = 3.759736666666667
$ws.Range("N8").Value = 11.27921
$ws.Range("O8").Value = 0.0683751702595819
$ws.Range("P8").Value = 0.06837517025958188
$ws.Range("Q8").Value = 28.77967882075334
$ws.Range("R8").Value = 259.01710938678
$ws.Range("S8").Value = 0.03333635049733082
$ws.Range("T8").Value = 0.03333635049733081

$ws.Range("G9").Value = 7.654706000000001
$ws.Range("I9").Value = 0.4875505299770593
$ws.Range("J9").Value = 0.4875505299770593
$ws.Range("O9").Value = 0.6514180024294648
$ws.Range("P9").Value = 0.6514180024294647
$ws.Range("Q9").Value = 274.1872644236594
$ws.Range("R9").Value = 2467.685379812935
$ws.Range("S9").Value = 0.3175991923210829
$ws.Range("T9").Value = 0.3175991923210829

$ws.Range("G10").Value = 7.654706000000001
$ws.Range("I10").Value = 0.4875505299770593
$ws.Range("J10").Value = 0.4875505299770593
$ws.Range("O10").Value = 0.2802068273109533
$ws.Range("P10").Value = 0.2802068273109533
$ws.Range("S10").Value = 0.1366149871586456
$ws.Range("T10").Value = 0.1366149871586456

